$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: new match data (07/08/2025 Godoy Cruz vs Gimnasia L.P.)
# Force text format on the date column so it is stored as a literal
# string (matching the source diff) instead of being auto-converted
# into a date serial number by Excel's smart input parsing.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "07/08/2025"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = "Godoy Cruz"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = "Gimnasia L.P."
$ws.Range("F21").Value = "W"
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 1.4
$ws.Range("L21").Value = 1.03
$ws.Range("M21").Value = 5
$ws.Range("N21").Value = 20
$ws.Range("O21").Value = 3
$ws.Range("P21").Value = 4
